# PacketIndex.xlsx edit: add a new "heading" row into the data/index table
# (row 11), pushing the existing fixQual/sat/End rows down by one, and
# update the dependent running-total / delta formulas + selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (just above the "fixQual" row), shifting the
# rows below it (fixQual, sat, End) down by one.
$ws.Rows.Item(11).Insert()

# Populate the new "heading" row with its index label, length/value (4),
# and the same running-sum / delta formulas used by the rest of the table.
$ws.Range("A11").Value = "heading"
$ws.Range("B11").Value = 4
$ws.Range("C11").Formula = "=SUM(B`$2:B11)"
$ws.Range("D11").Formula = "=C11-B11"

# Match the author's final selection (cell F13).
$ws.Range("F13").Select()
